# Rename the inline picture shapes (wp:docPr/@name) for the Pearson Edexcel
# logo in both footers and the BTEC logo in the header, matching the
# upstream commit that swapped image1.png <-> image2.png (Pearson logos)
# and image2.jpg -> image1.jpg (BTEC logo).
#
# NOTE: header/footer InlineShape collections must be walked paragraph by
# paragraph here - asking the Header/Footer Range directly for
# .InlineShapes(1) leaves the returned InlineShape's handle stale for a
# later property assignment in this runtime, so we re-seat it through the
# owning Paragraph's Range first.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlineShapesInRange($range, $newName) {
    $paras = $range.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        $shapes = $p.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            $shp.Name = $newName
        }
    }
}

# Footers(1) -> word/footer2.xml (docPr id="2", Pearson logo) : image1.png -> image2.png
Rename-InlineShapesInRange $sec.Footers(1).Range "image2.png"

# Footers(2) -> word/footer1.xml (docPr id="3", Pearson logo) : image1.png -> image2.png
Rename-InlineShapesInRange $sec.Footers(2).Range "image2.png"

# Headers(2) -> word/header1.xml (docPr id="1", BTEC logo) : image2.jpg -> image1.jpg
Rename-InlineShapesInRange $sec.Headers(2).Range "image1.jpg"
